# Update "Generate Report for Handback" timestamps across the three sheets.
# Sheet "Overview": G4 = Latest HO Xliff Generate Date for the 7332775b... row
# Sheet "zh-cn": H4 = Correspond Handoff Datetime, K4 = Correspond Handback DateTime (7332775b... row)
# Sheet "de-de": H4 = Correspond Handoff Datetime (shared text with Overview!G4), K4 = Correspond Handback DateTime

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (also shared with de-de!H4 text)
$wsOverview.Range("G4").Value = "2016-08-21 14:52:59"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-08-21 14:52:55"
$wsZhCn.Range("K4").Value = "2016-08-21 14:53:26"

# de-de sheet: Correspond Handoff Datetime (mirrors Overview!G4 value) / Correspond Handback DateTime
$wsDeDe.Range("H4").Value = "2016-08-21 14:52:59"
$wsDeDe.Range("K4").Value = "2016-08-21 14:53:32"
